# Edit script: reshuffle UK_mentalWellbeingByAgeGroup / UK_physicalWellbeingByAgeGroup /
# UK_lifeSatisfactionByAgeGroup sheets into UK_lifeSatisfactionByAgeGroup (refreshed values),
# UK_healthMCSByAgeGroup (new) and UK_healthPCSByAgeGroup (new), matching updated variable names
# and re-simulated validation statistics.

$wb = $excel.ActiveWorkbook

# --- Step 1: drop the two retired sheets ---
$wb.Worksheets.Item("UK_mentalWellbeingByAgeGroup").Delete()
$wb.Worksheets.Item("UK_physicalWellbeingByAgeGroup").Delete()

# --- Step 2: refresh UK_lifeSatisfactionByAgeGroup with newly-simulated values ---
$wsLS = $wb.Worksheets.Item("UK_lifeSatisfactionByAgeGroup")

$wsLS.Range("B2").Value = 5
$wsLS.Range("C2").Value = 5
$wsLS.Range("D2").Value = 4.9473670815003112
$wsLS.Range("E2").Value = 5.1040844640421081
$wsLS.Range("F2").Value = 5.0333282040384653
$wsLS.Range("G2").Value = 4.9030890916881713
$wsLS.Range("H2").Value = 4.9900372659891943
$wsLS.Range("I2").Value = 4.9532779135865628
$wsLS.Range("J2").Value = 5.0075712359057611
$wsLS.Range("K2").Value = 5.0862687149237997
$wsLS.Range("L2").Value = 4.9263923865663068
$wsLS.Range("M2").Value = 4.9954649101296829
$wsLS.Range("N2").Value = 5.1447935051748477
$wsLS.Range("O2").Value = 5.3392897410776232
$wsLS.Range("P2").Value = 5.564984983723992
$wsLS.Range("Q2").Value = 5.4003840827709606
$wsLS.Range("R2").Value = 5.6573521253079377
$wsLS.Range("S2").Value = 5.5216762931784897

$wsLS.Range("B3").Value = 5
$wsLS.Range("C3").Value = 5
$wsLS.Range("D3").Value = 4.9613212208943827
$wsLS.Range("E3").Value = 5.1095008113228877
$wsLS.Range("F3").Value = 5.0123925533412779
$wsLS.Range("G3").Value = 5.0301418339485693
$wsLS.Range("H3").Value = 5.0787445570875986
$wsLS.Range("I3").Value = 5.0471906266071578
$wsLS.Range("J3").Value = 4.9784911518287362
$wsLS.Range("K3").Value = 5.0527937191611896
$wsLS.Range("L3").Value = 4.9908818102540433
$wsLS.Range("M3").Value = 5.0700889811441261
$wsLS.Range("N3").Value = 5.1871266846892086
$wsLS.Range("O3").Value = 5.2664446748759062
$wsLS.Range("P3").Value = 5.4727455892632184
$wsLS.Range("Q3").Value = 5.5764862351364766
$wsLS.Range("R3").Value = 5.637319399278085
$wsLS.Range("S3").Value = 5.6706813648229328

$wsLS.Range("B4").Value = 5
$wsLS.Range("C4").Value = 5
$wsLS.Range("D4").Value = 4.9934442297436279
$wsLS.Range("E4").Value = 5.1199360069467383
$wsLS.Range("F4").Value = 5.0158787280590316
$wsLS.Range("G4").Value = 5.1020334320990806
$wsLS.Range("H4").Value = 5.1417166005043313
$wsLS.Range("I4").Value = 5.0629601510096327
$wsLS.Range("J4").Value = 4.9839173070781886
$wsLS.Range("K4").Value = 5.0653189584942906
$wsLS.Range("L4").Value = 4.9514677314738416
$wsLS.Range("M4").Value = 5.0755392121389864
$wsLS.Range("N4").Value = 5.2737492321617294
$wsLS.Range("O4").Value = 5.3052365966802766
$wsLS.Range("P4").Value = 5.5530234865234531
$wsLS.Range("Q4").Value = 5.6119841864645847
$wsLS.Range("R4").Value = 5.6535603285214906
$wsLS.Range("S4").Value = 5.659220700227019

$wsLS.Range("B5").Value = 5
$wsLS.Range("C5").Value = 5
$wsLS.Range("D5").Value = 4.9815058035494406
$wsLS.Range("E5").Value = 5.0800824937637499
$wsLS.Range("F5").Value = 5.0785856645403253
$wsLS.Range("G5").Value = 5.0306921075521798
$wsLS.Range("H5").Value = 5.0693314156142453
$wsLS.Range("I5").Value = 5.0260162923848952
$wsLS.Range("J5").Value = 4.9552083912161073
$wsLS.Range("K5").Value = 4.9866281814126356
$wsLS.Range("L5").Value = 4.9102675836177587
$wsLS.Range("M5").Value = 5.0076816238901456
$wsLS.Range("N5").Value = 5.2200884992539054
$wsLS.Range("O5").Value = 5.2579645682358072
$wsLS.Range("P5").Value = 5.4554878539910181
$wsLS.Range("Q5").Value = 5.5325465575076791
$wsLS.Range("R5").Value = 5.5180624594266359
$wsLS.Range("S5").Value = 5.4820374721529914

$wsLS.Range("B6").Value = 5
$wsLS.Range("C6").Value = 5
$wsLS.Range("D6").Value = 4.9654016268823069
$wsLS.Range("E6").Value = 5.0951873253109863
$wsLS.Range("F6").Value = 5.0597276811824976
$wsLS.Range("G6").Value = 5.0164525245646896
$wsLS.Range("H6").Value = 5.0498170502262143
$wsLS.Range("I6").Value = 5.0234953453249123
$wsLS.Range("J6").Value = 4.9786177855940954
$wsLS.Range("K6").Value = 5.0032311237537526
$wsLS.Range("L6").Value = 4.9474697827293479
$wsLS.Range("M6").Value = 5.010858783762683
$wsLS.Range("N6").Value = 5.2404017907079066
$wsLS.Range("O6").Value = 5.2816701538328106
$wsLS.Range("P6").Value = 5.4960577905424346
$wsLS.Range("Q6").Value = 5.5361494179142046
$wsLS.Range("R6").Value = 5.6349959629476958
$wsLS.Range("S6").Value = 5.5945908998065192

$wsLS.Range("B7").Value = 5
$wsLS.Range("C7").Value = 5
$wsLS.Range("D7").Value = 5.0181585209855486
$wsLS.Range("E7").Value = 5.1718471947531901
$wsLS.Range("F7").Value = 5.1596076851408856
$wsLS.Range("G7").Value = 5.1066017555767829
$wsLS.Range("H7").Value = 5.1037393310224877
$wsLS.Range("I7").Value = 5.1234941679568253
$wsLS.Range("J7").Value = 4.9873293532844221
$wsLS.Range("K7").Value = 5.0144566825019652
$wsLS.Range("L7").Value = 4.9925647355767024
$wsLS.Range("M7").Value = 5.0365709112652883
$wsLS.Range("N7").Value = 5.3058251598264494
$wsLS.Range("O7").Value = 5.3232571718954276
$wsLS.Range("P7").Value = 5.5378859501006898
$wsLS.Range("Q7").Value = 5.4951247475517997
$wsLS.Range("R7").Value = 5.5338439891927802
$wsLS.Range("S7").Value = 5.4576808937696262

$wsLS.Range("B8").Value = 5
$wsLS.Range("C8").Value = 5
$wsLS.Range("D8").Value = 5.091956714643378
$wsLS.Range("E8").Value = 5.1979045277818718
$wsLS.Range("F8").Value = 5.2042124879998202
$wsLS.Range("G8").Value = 5.1261608958290177
$wsLS.Range("H8").Value = 5.1757579794255868
$wsLS.Range("I8").Value = 5.1646315233332949
$wsLS.Range("J8").Value = 5.072366768255014
$wsLS.Range("K8").Value = 5.2336913369683877
$wsLS.Range("L8").Value = 5.0609537838680998
$wsLS.Range("M8").Value = 5.1514258032878733
$wsLS.Range("N8").Value = 5.3805427063853628
$wsLS.Range("O8").Value = 5.3578744924364248
$wsLS.Range("P8").Value = 5.5048935874882794
$wsLS.Range("Q8").Value = 5.5980584175023624
$wsLS.Range("R8").Value = 5.6329869947111044
$wsLS.Range("S8").Value = 5.525797385452389

$wsLS.Range("B9").Value = 5
$wsLS.Range("C9").Value = 5
$wsLS.Range("D9").Value = 5.1161207081428888
$wsLS.Range("E9").Value = 5.1586791710112356
$wsLS.Range("F9").Value = 5.2322350720679003
$wsLS.Range("G9").Value = 5.0832157211609754
$wsLS.Range("H9").Value = 5.2577953512742202
$wsLS.Range("I9").Value = 5.156980686444613
$wsLS.Range("J9").Value = 5.1187839229182117
$wsLS.Range("K9").Value = 5.2070515266428181
$wsLS.Range("L9").Value = 5.0712504208004123
$wsLS.Range("M9").Value = 5.205072319555919
$wsLS.Range("N9").Value = 5.4061642306260991
$wsLS.Range("O9").Value = 5.4122209526646357
$wsLS.Range("P9").Value = 5.6075850567044698
$wsLS.Range("Q9").Value = 5.5589451975493267
$wsLS.Range("R9").Value = 5.6058398265117271
$wsLS.Range("S9").Value = 5.5649840715378822

$wsLS.Range("B10").Value = 5
$wsLS.Range("C10").Value = 5
$wsLS.Range("D10").Value = 5.0764589583984776
$wsLS.Range("E10").Value = 5.1371618670971477
$wsLS.Range("F10").Value = 5.1415114541260714
$wsLS.Range("G10").Value = 5.0393908872109749
$wsLS.Range("H10").Value = 5.1301997611431061
$wsLS.Range("I10").Value = 5.1695715747398694
$wsLS.Range("J10").Value = 5.0132614581149566
$wsLS.Range("K10").Value = 5.0950120107390138
$wsLS.Range("L10").Value = 5.0081136352843769
$wsLS.Range("M10").Value = 5.0053232542660089
$wsLS.Range("N10").Value = 5.2406530729039131
$wsLS.Range("O10").Value = 5.2637436520297092
$wsLS.Range("P10").Value = 5.4207346222576103
$wsLS.Range("Q10").Value = 5.366055383274233
$wsLS.Range("R10").Value = 5.412244371410142
$wsLS.Range("S10").Value = 5.3903643381577249

$wsLS.Range("B11").Value = 5
$wsLS.Range("C11").Value = 5
$wsLS.Range("D11").Value = 5.0269010706692763
$wsLS.Range("E11").Value = 5.1348274908599354
$wsLS.Range("F11").Value = 5.0735101031514391
$wsLS.Range("G11").Value = 5.013711248729372
$wsLS.Range("H11").Value = 5.0312839818835444
$wsLS.Range("I11").Value = 5.0187757957899093
$wsLS.Range("J11").Value = 4.8780584621308751
$wsLS.Range("K11").Value = 4.912059750537578
$wsLS.Range("L11").Value = 4.782507941384158
$wsLS.Range("M11").Value = 4.8404125767247228
$wsLS.Range("N11").Value = 5.2115792323504992
$wsLS.Range("O11").Value = 5.231789563674754
$wsLS.Range("P11").Value = 5.2659253126512509
$wsLS.Range("Q11").Value = 5.2849310471579436
$wsLS.Range("R11").Value = 5.2209533631997944
$wsLS.Range("S11").Value = 5.274529490219324

$wsLS.Range("B12").Value = 5
$wsLS.Range("C12").Value = 5
$wsLS.Range("D12").Value = 5.0308304616237063
$wsLS.Range("E12").Value = 5.1090067791443943
$wsLS.Range("F12").Value = 5.0794961760364954
$wsLS.Range("G12").Value = 5.0289939318549708
$wsLS.Range("H12").Value = 5.1014296287027738
$wsLS.Range("I12").Value = 5.0318479219292476
$wsLS.Range("J12").Value = 4.9400838483586043
$wsLS.Range("K12").Value = 4.9464984353218124
$wsLS.Range("L12").Value = 4.911621954964656
$wsLS.Range("M12").Value = 4.9646400874541063
$wsLS.Range("N12").Value = 5.3228031139371224
$wsLS.Range("O12").Value = 5.2921515435831052
$wsLS.Range("P12").Value = 5.3400803544434332
$wsLS.Range("Q12").Value = 5.3663605222657997
$wsLS.Range("R12").Value = 5.2207790839758994
$wsLS.Range("S12").Value = 5.3171145739591719

$wsLS.Range("B13").Value = 5
$wsLS.Range("C13").Value = 5
$wsLS.Range("D13").Value = 5.1176851247528212
$wsLS.Range("E13").Value = 5.1671432816770348
$wsLS.Range("F13").Value = 5.1569093627650764
$wsLS.Range("G13").Value = 5.1183831664895942
$wsLS.Range("H13").Value = 5.0922362798196943
$wsLS.Range("I13").Value = 5.0810507047778053
$wsLS.Range("J13").Value = 4.9971914506053681
$wsLS.Range("K13").Value = 5.0848699512533857
$wsLS.Range("L13").Value = 5.0091007177902966
$wsLS.Range("M13").Value = 5.0797028974482403
$wsLS.Range("N13").Value = 5.3715733788656959
$wsLS.Range("O13").Value = 5.343952981175339
$wsLS.Range("P13").Value = 5.4784648453319207
$wsLS.Range("Q13").Value = 5.4776614581525056
$wsLS.Range("R13").Value = 5.3830223159770894
$wsLS.Range("S13").Value = 5.4919357866675567

# --- Step 3: create UK_healthMCSByAgeGroup as a new sheet after UK_lifeSatisfactionByAgeGroup ---
$wsMCS = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsLS)
$wsMCS.Name = "UK_healthMCSByAgeGroup"
$wsMCS.Range("A1").Value = "Year"
$wsMCS.Range("B1").Value = "health_mcs_score_female_0_9"
$wsMCS.Range("C1").Value = "health_mcs_score_male_0_9"
$wsMCS.Range("D1").Value = "health_mcs_score_female_10_19"
$wsMCS.Range("E1").Value = "health_mcs_score_male_10_19"
$wsMCS.Range("F1").Value = "health_mcs_score_female_20_29"
$wsMCS.Range("G1").Value = "health_mcs_score_male_20_29"
$wsMCS.Range("H1").Value = "health_mcs_score_female_30_39"
$wsMCS.Range("I1").Value = "health_mcs_score_male_30_39"
$wsMCS.Range("J1").Value = "health_mcs_score_female_40_49"
$wsMCS.Range("K1").Value = "health_mcs_score_male_40_49"
$wsMCS.Range("L1").Value = "health_mcs_score_female_50_59"
$wsMCS.Range("M1").Value = "health_mcs_score_male_50_59"
$wsMCS.Range("N1").Value = "health_mcs_score_female_60_69"
$wsMCS.Range("O1").Value = "health_mcs_score_male_60_69"
$wsMCS.Range("P1").Value = "health_mcs_score_female_70_79"
$wsMCS.Range("Q1").Value = "health_mcs_score_male_70_79"
$wsMCS.Range("R1").Value = "health_mcs_score_female_80_100"
$wsMCS.Range("S1").Value = "health_mcs_score_male_80_100"

$wsMCS.Range("A2").Value = 2022
$wsMCS.Range("B2").Value = 42.387932395771813
$wsMCS.Range("C2").Value = 44.138288521782243
$wsMCS.Range("D2").Value = 41.773307197990981
$wsMCS.Range("E2").Value = 46.542943401439921
$wsMCS.Range("F2").Value = 42.777534136271512
$wsMCS.Range("G2").Value = 45.389405052402829
$wsMCS.Range("H2").Value = 42.801475275462103
$wsMCS.Range("I2").Value = 44.521091951445648
$wsMCS.Range("J2").Value = 44.823522520167373
$wsMCS.Range("K2").Value = 47.541100806181348
$wsMCS.Range("L2").Value = 46.129113843686198
$wsMCS.Range("M2").Value = 48.322070178801177
$wsMCS.Range("N2").Value = 48.232054546044942
$wsMCS.Range("O2").Value = 50.792932958267123
$wsMCS.Range("P2").Value = 50.347353923552127
$wsMCS.Range("Q2").Value = 51.371380455976599
$wsMCS.Range("R2").Value = 52.547699544956672
$wsMCS.Range("S2").Value = 52.209875546611748

$wsMCS.Range("A3").Value = 2021
$wsMCS.Range("B3").Value = 42.464496033292257
$wsMCS.Range("C3").Value = 44.225424703842492
$wsMCS.Range("D3").Value = 41.502513861235428
$wsMCS.Range("E3").Value = 45.820091774316211
$wsMCS.Range("F3").Value = 40.829191259345258
$wsMCS.Range("G3").Value = 45.226492817854457
$wsMCS.Range("H3").Value = 43.081382047195113
$wsMCS.Range("I3").Value = 45.626545613685067
$wsMCS.Range("J3").Value = 44.675078684324831
$wsMCS.Range("K3").Value = 47.128095649974242
$wsMCS.Range("L3").Value = 46.531513209736858
$wsMCS.Range("M3").Value = 48.790862004833748
$wsMCS.Range("N3").Value = 48.466301065923417
$wsMCS.Range("O3").Value = 50.670435496496793
$wsMCS.Range("P3").Value = 50.229495650873609
$wsMCS.Range("Q3").Value = 52.4129361430739
$wsMCS.Range("R3").Value = 51.560535747375177
$wsMCS.Range("S3").Value = 52.971874522853781

$wsMCS.Range("A4").Value = 2020
$wsMCS.Range("B4").Value = 42.664327248462897
$wsMCS.Range("C4").Value = 44.54530229454452
$wsMCS.Range("D4").Value = 41.807293336803973
$wsMCS.Range("E4").Value = 45.70625623649115
$wsMCS.Range("F4").Value = 41.420042049533713
$wsMCS.Range("G4").Value = 46.222960821502632
$wsMCS.Range("H4").Value = 43.092662753102132
$wsMCS.Range("I4").Value = 45.494703092555277
$wsMCS.Range("J4").Value = 44.639196202321664
$wsMCS.Range("K4").Value = 47.779246242625788
$wsMCS.Range("L4").Value = 46.091694328398653
$wsMCS.Range("M4").Value = 48.839218166202699
$wsMCS.Range("N4").Value = 48.523156481090467
$wsMCS.Range("O4").Value = 50.550717763173807
$wsMCS.Range("P4").Value = 50.37686646006577
$wsMCS.Range("Q4").Value = 52.596218539254231
$wsMCS.Range("R4").Value = 51.528419025308743
$wsMCS.Range("S4").Value = 52.607787299207757

$wsMCS.Range("A5").Value = 2019
$wsMCS.Range("B5").Value = 43.006434440124359
$wsMCS.Range("C5").Value = 45.047069465022062
$wsMCS.Range("D5").Value = 43.058389214666398
$wsMCS.Range("E5").Value = 46.407244269732537
$wsMCS.Range("F5").Value = 42.574386139931292
$wsMCS.Range("G5").Value = 45.896448302684988
$wsMCS.Range("H5").Value = 44.115234792169574
$wsMCS.Range("I5").Value = 46.258173661457427
$wsMCS.Range("J5").Value = 45.964596037775287
$wsMCS.Range("K5").Value = 47.645333075612527
$wsMCS.Range("L5").Value = 47.048959422028908
$wsMCS.Range("M5").Value = 48.816186699032208
$wsMCS.Range("N5").Value = 49.31153919037672
$wsMCS.Range("O5").Value = 51.040581438438643
$wsMCS.Range("P5").Value = 50.980881872742877
$wsMCS.Range("Q5").Value = 52.598341383166087
$wsMCS.Range("R5").Value = 51.121015599843822
$wsMCS.Range("S5").Value = 52.270196693878312

$wsMCS.Range("A6").Value = 2018
$wsMCS.Range("B6").Value = 43.529125662451079
$wsMCS.Range("C6").Value = 45.59834374594341
$wsMCS.Range("D6").Value = 43.051482483468327
$wsMCS.Range("E6").Value = 46.533378572133067
$wsMCS.Range("F6").Value = 43.329103229042254
$wsMCS.Range("G6").Value = 46.717118204810383
$wsMCS.Range("H6").Value = 44.469032141982808
$wsMCS.Range("I6").Value = 47.045121607523079
$wsMCS.Range("J6").Value = 46.249414344467382
$wsMCS.Range("K6").Value = 48.147328468818387
$wsMCS.Range("L6").Value = 47.264286837860098
$wsMCS.Range("M6").Value = 49.248939040415827
$wsMCS.Range("N6").Value = 49.611796012418637
$wsMCS.Range("O6").Value = 51.519424742256867
$wsMCS.Range("P6").Value = 51.231645125806502
$wsMCS.Range("Q6").Value = 52.843593660276568
$wsMCS.Range("R6").Value = 51.668324186088249
$wsMCS.Range("S6").Value = 52.2509479274391

$wsMCS.Range("A7").Value = 2017
$wsMCS.Range("B7").Value = 43.937505160864127
$wsMCS.Range("C7").Value = 45.94402125910198
$wsMCS.Range("D7").Value = 44.056619984096621
$wsMCS.Range("E7").Value = 47.49306860895809
$wsMCS.Range("F7").Value = 44.848755037985832
$wsMCS.Range("G7").Value = 47.250134016471158
$wsMCS.Range("H7").Value = 45.527634810511238
$wsMCS.Range("I7").Value = 47.307872058276708
$wsMCS.Range("J7").Value = 46.4804220920212
$wsMCS.Range("K7").Value = 48.528704182879359
$wsMCS.Range("L7").Value = 47.686701428742538
$wsMCS.Range("M7").Value = 49.764110630171068
$wsMCS.Range("N7").Value = 50.097950226771012
$wsMCS.Range("O7").Value = 51.622395651330933
$wsMCS.Range("P7").Value = 51.3655028210338
$wsMCS.Range("Q7").Value = 52.65759678376844
$wsMCS.Range("R7").Value = 51.558809788628842
$wsMCS.Range("S7").Value = 52.871498867811447

$wsMCS.Range("A8").Value = 2016
$wsMCS.Range("B8").Value = 44.147445654864747
$wsMCS.Range("C8").Value = 46.029077244037452
$wsMCS.Range("D8").Value = 44.656633106362932
$wsMCS.Range("E8").Value = 48.037692329313643
$wsMCS.Range("F8").Value = 45.535307215444753
$wsMCS.Range("G8").Value = 47.812651970088297
$wsMCS.Range("H8").Value = 46.032549647892452
$wsMCS.Range("I8").Value = 48.223330696476609
$wsMCS.Range("J8").Value = 47.060051959350957
$wsMCS.Range("K8").Value = 49.284716792966861
$wsMCS.Range("L8").Value = 48.048785001043747
$wsMCS.Range("M8").Value = 50.215713067347707
$wsMCS.Range("N8").Value = 50.592566764768883
$wsMCS.Range("O8").Value = 51.894057739941317
$wsMCS.Range("P8").Value = 51.29564878630535
$wsMCS.Range("Q8").Value = 52.906243023153188
$wsMCS.Range("R8").Value = 51.442533076120291
$wsMCS.Range("S8").Value = 52.317513309095908

$wsMCS.Range("A9").Value = 2015
$wsMCS.Range("B9").Value = 44.439814869756127
$wsMCS.Range("C9").Value = 46.297600937837203
$wsMCS.Range("D9").Value = 45.366196851037273
$wsMCS.Range("E9").Value = 48.363838392564688
$wsMCS.Range("F9").Value = 46.203452804171043
$wsMCS.Range("G9").Value = 48.354172831788482
$wsMCS.Range("H9").Value = 46.653392208242813
$wsMCS.Range("I9").Value = 48.780592775721637
$wsMCS.Range("J9").Value = 47.505806676354538
$wsMCS.Range("K9").Value = 49.854920056439937
$wsMCS.Range("L9").Value = 48.053648285114548
$wsMCS.Range("M9").Value = 50.858326168693459
$wsMCS.Range("N9").Value = 50.734927396152912
$wsMCS.Range("O9").Value = 52.473127122016407
$wsMCS.Range("P9").Value = 51.905041049588633
$wsMCS.Range("Q9").Value = 53.512632035519587
$wsMCS.Range("R9").Value = 51.883297169569147
$wsMCS.Range("S9").Value = 52.835801726962899

$wsMCS.Range("A10").Value = 2014
$wsMCS.Range("B10").Value = 44.681266024181546
$wsMCS.Range("C10").Value = 46.560294737720042
$wsMCS.Range("D10").Value = 45.583670127318257
$wsMCS.Range("E10").Value = 48.157520392592041
$wsMCS.Range("F10").Value = 46.287988042391802
$wsMCS.Range("G10").Value = 48.746008167841723
$wsMCS.Range("H10").Value = 47.231488748230028
$wsMCS.Range("I10").Value = 49.205504873377137
$wsMCS.Range("J10").Value = 48.06338307526741
$wsMCS.Range("K10").Value = 50.217244317996602
$wsMCS.Range("L10").Value = 48.425787893396922
$wsMCS.Range("M10").Value = 50.66533178265005
$wsMCS.Range("N10").Value = 50.984123961689512
$wsMCS.Range("O10").Value = 52.671143169602168
$wsMCS.Range("P10").Value = 52.059981177263737
$wsMCS.Range("Q10").Value = 53.378006501791127
$wsMCS.Range("R10").Value = 52.068218410787189
$wsMCS.Range("S10").Value = 52.859142460492649

$wsMCS.Range("A11").Value = 2013
$wsMCS.Range("B11").Value = 44.542020120314547
$wsMCS.Range("C11").Value = 46.530975791210373
$wsMCS.Range("D11").Value = 45.474999507306428
$wsMCS.Range("E11").Value = 48.116704974877038
$wsMCS.Range("F11").Value = 46.769707805973567
$wsMCS.Range("G11").Value = 48.859659902286559
$wsMCS.Range("H11").Value = 46.941636462130333
$wsMCS.Range("I11").Value = 49.3036768597471
$wsMCS.Range("J11").Value = 47.680287420579518
$wsMCS.Range("K11").Value = 49.763054201837072
$wsMCS.Range("L11").Value = 48.166320573784937
$wsMCS.Range("M11").Value = 50.470337390291071
$wsMCS.Range("N11").Value = 51.053922971664498
$wsMCS.Range("O11").Value = 52.710970851844358
$wsMCS.Range("P11").Value = 51.550802084984781
$wsMCS.Range("Q11").Value = 53.08815107068704
$wsMCS.Range("R11").Value = 51.74764864150584
$wsMCS.Range("S11").Value = 53.680865849140467

$wsMCS.Range("A12").Value = 2012
$wsMCS.Range("B12").Value = 44.591236272660048
$wsMCS.Range("C12").Value = 46.644496096770759
$wsMCS.Range("D12").Value = 45.69556029409825
$wsMCS.Range("E12").Value = 48.027714825779938
$wsMCS.Range("F12").Value = 46.346896048571068
$wsMCS.Range("G12").Value = 48.964676610131711
$wsMCS.Range("H12").Value = 47.316753700880533
$wsMCS.Range("I12").Value = 49.646826377060798
$wsMCS.Range("J12").Value = 47.74657483065922
$wsMCS.Range("K12").Value = 49.926951069883778
$wsMCS.Range("L12").Value = 48.28048046816204
$wsMCS.Range("M12").Value = 50.94644357265252
$wsMCS.Range("N12").Value = 51.256791372488621
$wsMCS.Range("O12").Value = 53.06433202168656
$wsMCS.Range("P12").Value = 51.669341662263001
$wsMCS.Range("Q12").Value = 53.517222965794552
$wsMCS.Range("R12").Value = 51.562034680086761
$wsMCS.Range("S12").Value = 53.278209520200839

$wsMCS.Range("A13").Value = 2011
$wsMCS.Range("B13").Value = 44.961598280175913
$wsMCS.Range("C13").Value = 46.929734834477777
$wsMCS.Range("D13").Value = 46.246464861618051
$wsMCS.Range("E13").Value = 48.508359773925093
$wsMCS.Range("F13").Value = 47.381383476146119
$wsMCS.Range("G13").Value = 49.612984607556193
$wsMCS.Range("H13").Value = 47.781103558979012
$wsMCS.Range("I13").Value = 49.748321519986668
$wsMCS.Range("J13").Value = 48.493845970320727
$wsMCS.Range("K13").Value = 50.381672934926883
$wsMCS.Range("L13").Value = 48.701679532754753
$wsMCS.Range("M13").Value = 50.987328315529247
$wsMCS.Range("N13").Value = 51.295903211173098
$wsMCS.Range("O13").Value = 52.735782735604822
$wsMCS.Range("P13").Value = 51.528728700040439
$wsMCS.Range("Q13").Value = 53.291376624384448
$wsMCS.Range("R13").Value = 51.948619429818677
$wsMCS.Range("S13").Value = 53.622609313269479

# --- Step 4: create UK_healthPCSByAgeGroup as a new sheet after UK_healthMCSByAgeGroup ---
$wsPCS = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsMCS)
$wsPCS.Name = "UK_healthPCSByAgeGroup"
$wsPCS.Range("A1").Value = "Year"
$wsPCS.Range("B1").Value = "health_pcs_score_female_0_9"
$wsPCS.Range("C1").Value = "health_pcs_score_male_0_9"
$wsPCS.Range("D1").Value = "health_pcs_score_female_10_19"
$wsPCS.Range("E1").Value = "health_pcs_score_male_10_19"
$wsPCS.Range("F1").Value = "health_pcs_score_female_20_29"
$wsPCS.Range("G1").Value = "health_pcs_score_male_20_29"
$wsPCS.Range("H1").Value = "health_pcs_score_female_30_39"
$wsPCS.Range("I1").Value = "health_pcs_score_male_30_39"
$wsPCS.Range("J1").Value = "health_pcs_score_female_40_49"
$wsPCS.Range("K1").Value = "health_pcs_score_male_40_49"
$wsPCS.Range("L1").Value = "health_pcs_score_female_50_59"
$wsPCS.Range("M1").Value = "health_pcs_score_male_50_59"
$wsPCS.Range("N1").Value = "health_pcs_score_female_60_69"
$wsPCS.Range("O1").Value = "health_pcs_score_male_60_69"
$wsPCS.Range("P1").Value = "health_pcs_score_female_70_79"
$wsPCS.Range("Q1").Value = "health_pcs_score_male_70_79"
$wsPCS.Range("R1").Value = "health_pcs_score_female_80_100"
$wsPCS.Range("S1").Value = "health_pcs_score_male_80_100"

$wsPCS.Range("A2").Value = 2022
$wsPCS.Range("B2").Value = 61.321330893104722
$wsPCS.Range("C2").Value = 61.867377507998633
$wsPCS.Range("D2").Value = 56.764586675448697
$wsPCS.Range("E2").Value = 57.880819762251633
$wsPCS.Range("F2").Value = 53.984943883071587
$wsPCS.Range("G2").Value = 54.047871288023558
$wsPCS.Range("H2").Value = 52.138438618660018
$wsPCS.Range("I2").Value = 54.058104611712473
$wsPCS.Range("J2").Value = 49.586578727808757
$wsPCS.Range("K2").Value = 51.540627438195791
$wsPCS.Range("L2").Value = 47.673338993140653
$wsPCS.Range("M2").Value = 49.951567523232093
$wsPCS.Range("N2").Value = 45.927564889526089
$wsPCS.Range("O2").Value = 46.966873931872037
$wsPCS.Range("P2").Value = 44.300520229199797
$wsPCS.Range("Q2").Value = 45.473690233874208
$wsPCS.Range("R2").Value = 36.21869549469784
$wsPCS.Range("S2").Value = 40.600546266142317

$wsPCS.Range("A3").Value = 2021
$wsPCS.Range("B3").Value = 61.357144096633377
$wsPCS.Range("C3").Value = 61.902901058911553
$wsPCS.Range("D3").Value = 57.133249886669127
$wsPCS.Range("E3").Value = 58.298727631147919
$wsPCS.Range("F3").Value = 54.86131107074052
$wsPCS.Range("G3").Value = 55.313257724363012
$wsPCS.Range("H3").Value = 52.791064121641469
$wsPCS.Range("I3").Value = 54.132243653702332
$wsPCS.Range("J3").Value = 51.187332885673847
$wsPCS.Range("K3").Value = 52.222830827730739
$wsPCS.Range("L3").Value = 48.547236445764852
$wsPCS.Range("M3").Value = 49.867040256528419
$wsPCS.Range("N3").Value = 46.639964365263431
$wsPCS.Range("O3").Value = 47.257184952541103
$wsPCS.Range("P3").Value = 43.904480006020847
$wsPCS.Range("Q3").Value = 44.943348444776582
$wsPCS.Range("R3").Value = 36.722366022520013
$wsPCS.Range("S3").Value = 41.054562072107608

$wsPCS.Range("A4").Value = 2020
$wsPCS.Range("B4").Value = 61.341869578428103
$wsPCS.Range("C4").Value = 61.800356131297121
$wsPCS.Range("D4").Value = 57.616410266410853
$wsPCS.Range("E4").Value = 58.036695513879359
$wsPCS.Range("F4").Value = 54.661116613150483
$wsPCS.Range("G4").Value = 55.141139837533188
$wsPCS.Range("H4").Value = 53.432812061518433
$wsPCS.Range("I4").Value = 54.615179989521003
$wsPCS.Range("J4").Value = 51.066950654443268
$wsPCS.Range("K4").Value = 52.452657914254537
$wsPCS.Range("L4").Value = 48.485171909284418
$wsPCS.Range("M4").Value = 49.993935882030407
$wsPCS.Range("N4").Value = 46.651598263447632
$wsPCS.Range("O4").Value = 46.948861251379093
$wsPCS.Range("P4").Value = 43.263373054451137
$wsPCS.Range("Q4").Value = 44.892190777145153
$wsPCS.Range("R4").Value = 37.387707842379463
$wsPCS.Range("S4").Value = 39.843760440890513

$wsPCS.Range("A5").Value = 2019
$wsPCS.Range("B5").Value = 60.941166711149897
$wsPCS.Range("C5").Value = 61.473966520686282
$wsPCS.Range("D5").Value = 57.007293502740957
$wsPCS.Range("E5").Value = 57.357188497609933
$wsPCS.Range("F5").Value = 54.551171899142638
$wsPCS.Range("G5").Value = 55.376297040667417
$wsPCS.Range("H5").Value = 52.110855684300297
$wsPCS.Range("I5").Value = 53.988717340536468
$wsPCS.Range("J5").Value = 50.75825033935849
$wsPCS.Range("K5").Value = 52.143815385104993
$wsPCS.Range("L5").Value = 48.494789468413792
$wsPCS.Range("M5").Value = 49.727577806226513
$wsPCS.Range("N5").Value = 46.06531349099977
$wsPCS.Range("O5").Value = 46.636477672866441
$wsPCS.Range("P5").Value = 43.03746812475385
$wsPCS.Range("Q5").Value = 45.386180889894227
$wsPCS.Range("R5").Value = 38.426777416205297
$wsPCS.Range("S5").Value = 40.319883172604648

$wsPCS.Range("A6").Value = 2018
$wsPCS.Range("B6").Value = 60.635365785225559
$wsPCS.Range("C6").Value = 61.248346908004009
$wsPCS.Range("D6").Value = 56.911311827093421
$wsPCS.Range("E6").Value = 57.274800716973843
$wsPCS.Range("F6").Value = 54.500671160666307
$wsPCS.Range("G6").Value = 55.39140846613649
$wsPCS.Range("H6").Value = 52.743263926108611
$wsPCS.Range("I6").Value = 54.389269693042003
$wsPCS.Range("J6").Value = 50.735595742030547
$wsPCS.Range("K6").Value = 52.271001932054723
$wsPCS.Range("L6").Value = 48.404105506162317
$wsPCS.Range("M6").Value = 49.86005623210206
$wsPCS.Range("N6").Value = 45.740170211708239
$wsPCS.Range("O6").Value = 47.17231879925842
$wsPCS.Range("P6").Value = 43.51141487979892
$wsPCS.Range("Q6").Value = 44.997207449490631
$wsPCS.Range("R6").Value = 38.157000646155517
$wsPCS.Range("S6").Value = 40.266696849462868

$wsPCS.Range("A7").Value = 2017
$wsPCS.Range("B7").Value = 60.629752510760397
$wsPCS.Range("C7").Value = 61.125522940715257
$wsPCS.Range("D7").Value = 56.690444221483389
$wsPCS.Range("E7").Value = 57.354484872420883
$wsPCS.Range("F7").Value = 54.427171225602542
$wsPCS.Range("G7").Value = 55.513559482436293
$wsPCS.Range("H7").Value = 52.934495102443798
$wsPCS.Range("I7").Value = 54.375877321121713
$wsPCS.Range("J7").Value = 50.887068401880789
$wsPCS.Range("K7").Value = 52.123470590203937
$wsPCS.Range("L7").Value = 48.494663904676877
$wsPCS.Range("M7").Value = 49.213821455170859
$wsPCS.Range("N7").Value = 45.925983653848668
$wsPCS.Range("O7").Value = 46.761543468884803
$wsPCS.Range("P7").Value = 43.097040644837953
$wsPCS.Range("Q7").Value = 44.723569855456873
$wsPCS.Range("R7").Value = 37.488517488035733
$wsPCS.Range("S7").Value = 39.418748960193533

$wsPCS.Range("A8").Value = 2016
$wsPCS.Range("B8").Value = 60.353619475508992
$wsPCS.Range("C8").Value = 60.954615467126899
$wsPCS.Range("D8").Value = 57.094903860118123
$wsPCS.Range("E8").Value = 57.564811685539411
$wsPCS.Range("F8").Value = 54.333622921988479
$wsPCS.Range("G8").Value = 55.494772969568302
$wsPCS.Range("H8").Value = 53.077178152767203
$wsPCS.Range("I8").Value = 54.48238978612595
$wsPCS.Range("J8").Value = 50.99080570198916
$wsPCS.Range("K8").Value = 52.751094072235837
$wsPCS.Range("L8").Value = 48.542738083565709
$wsPCS.Range("M8").Value = 49.605959955612121
$wsPCS.Range("N8").Value = 46.142610622430787
$wsPCS.Range("O8").Value = 47.128155291476851
$wsPCS.Range("P8").Value = 42.450699124394859
$wsPCS.Range("Q8").Value = 44.397583991917188
$wsPCS.Range("R8").Value = 37.524424307811387
$wsPCS.Range("S8").Value = 39.84915810191756

$wsPCS.Range("A9").Value = 2015
$wsPCS.Range("B9").Value = 60.138571611488821
$wsPCS.Range("C9").Value = 60.642514057769567
$wsPCS.Range("D9").Value = 56.89555439361969
$wsPCS.Range("E9").Value = 57.597302503673397
$wsPCS.Range("F9").Value = 54.601265437888003
$wsPCS.Range("G9").Value = 55.858637125374123
$wsPCS.Range("H9").Value = 53.070250358958369
$wsPCS.Range("I9").Value = 54.321014139721989
$wsPCS.Range("J9").Value = 51.197581475806501
$wsPCS.Range("K9").Value = 52.645486514944231
$wsPCS.Range("L9").Value = 48.720476826273753
$wsPCS.Range("M9").Value = 50.198515086826717
$wsPCS.Range("N9").Value = 46.091826514558477
$wsPCS.Range("O9").Value = 47.225803882053803
$wsPCS.Range("P9").Value = 42.994746282706942
$wsPCS.Range("Q9").Value = 44.304782507942598
$wsPCS.Range("R9").Value = 37.893250863651851
$wsPCS.Range("S9").Value = 39.352172199681789

$wsPCS.Range("A10").Value = 2014
$wsPCS.Range("B10").Value = 59.831606697301211
$wsPCS.Range("C10").Value = 60.40838619398663
$wsPCS.Range("D10").Value = 56.547492209419481
$wsPCS.Range("E10").Value = 57.10837087074983
$wsPCS.Range("F10").Value = 54.334526745312203
$wsPCS.Range("G10").Value = 55.555292223645409
$wsPCS.Range("H10").Value = 53.149158236566358
$wsPCS.Range("I10").Value = 54.257084619388003
$wsPCS.Range("J10").Value = 51.353726772434847
$wsPCS.Range("K10").Value = 52.632467897088063
$wsPCS.Range("L10").Value = 48.774476384309317
$wsPCS.Range("M10").Value = 50.066181742886549
$wsPCS.Range("N10").Value = 45.950606838335823
$wsPCS.Range("O10").Value = 46.728700114406287
$wsPCS.Range("P10").Value = 42.636977887616432
$wsPCS.Range("Q10").Value = 43.902805168742383
$wsPCS.Range("R10").Value = 37.365107853917713
$wsPCS.Range("S10").Value = 39.518952871319449

$wsPCS.Range("A11").Value = 2013
$wsPCS.Range("B11").Value = 59.783508447303511
$wsPCS.Range("C11").Value = 60.436915488348554
$wsPCS.Range("D11").Value = 56.829776664471439
$wsPCS.Range("E11").Value = 57.150289528515103
$wsPCS.Range("F11").Value = 54.337327896295903
$wsPCS.Range("G11").Value = 55.235208975054903
$wsPCS.Range("H11").Value = 53.472931022009817
$wsPCS.Range("I11").Value = 54.317610288912768
$wsPCS.Range("J11").Value = 51.485480214747447
$wsPCS.Range("K11").Value = 52.754557462216731
$wsPCS.Range("L11").Value = 49.00095196342582
$wsPCS.Range("M11").Value = 50.006249458533937
$wsPCS.Range("N11").Value = 46.058320867747263
$wsPCS.Range("O11").Value = 47.23260953457747
$wsPCS.Range("P11").Value = 42.413534290768773
$wsPCS.Range("Q11").Value = 43.548990647552571
$wsPCS.Range("R11").Value = 37.631706410545661
$wsPCS.Range("S11").Value = 40.401702448873728

$wsPCS.Range("A12").Value = 2012
$wsPCS.Range("B12").Value = 59.720111431329322
$wsPCS.Range("C12").Value = 60.426153372133207
$wsPCS.Range("D12").Value = 56.287818002928773
$wsPCS.Range("E12").Value = 57.02553251703246
$wsPCS.Range("F12").Value = 54.128384261371188
$wsPCS.Range("G12").Value = 55.175832994442217
$wsPCS.Range("H12").Value = 53.273057494314664
$wsPCS.Range("I12").Value = 54.078994391401977
$wsPCS.Range("J12").Value = 51.678295107772357
$wsPCS.Range("K12").Value = 52.683593648897279
$wsPCS.Range("L12").Value = 49.02593647569568
$wsPCS.Range("M12").Value = 50.358925421940448
$wsPCS.Range("N12").Value = 46.343569742663597
$wsPCS.Range("O12").Value = 46.997701802793571
$wsPCS.Range("P12").Value = 42.311555511970496
$wsPCS.Range("Q12").Value = 44.079763540684063
$wsPCS.Range("R12").Value = 38.154248870243968
$wsPCS.Range("S12").Value = 40.839599545303862

$wsPCS.Range("A13").Value = 2011
$wsPCS.Range("B13").Value = 59.668692788581943
$wsPCS.Range("C13").Value = 60.491844780602698
$wsPCS.Range("D13").Value = 56.433131239288443
$wsPCS.Range("E13").Value = 57.157544062819973
$wsPCS.Range("F13").Value = 54.432116597061267
$wsPCS.Range("G13").Value = 55.388803205042123
$wsPCS.Range("H13").Value = 53.513171015649952
$wsPCS.Range("I13").Value = 54.44162951266496
$wsPCS.Range("J13").Value = 51.604394488828213
$wsPCS.Range("K13").Value = 52.692942333289437
$wsPCS.Range("L13").Value = 48.823542187215999
$wsPCS.Range("M13").Value = 49.982755911747191
$wsPCS.Range("N13").Value = 46.147637881215289
$wsPCS.Range("O13").Value = 46.895231410894162
$wsPCS.Range("P13").Value = 42.085259168606242
$wsPCS.Range("Q13").Value = 43.141845370236133
$wsPCS.Range("R13").Value = 37.575338617385214
$wsPCS.Range("S13").Value = 40.297658061501579

# --- Step 5: make the last sheet (UK_healthPCSByAgeGroup) the active / selected tab ---
$wsPCS.Activate()

Write-Host "Done. Sheets now:"
foreach ($s in $wb.Worksheets) { Write-Host $s.Name }
